# Map table name (column B) to primary/composite key (column C).
# Order matters: it controls the order new entries land in the shared
# string table, which must match addr -> value pairs below exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orderedEdits = @(
    @("C3",  "actor_id"),
    @("C4",  "address_id"),
    @("C5",  "category_id"),
    @("C6",  "city_id"),
    @("C7",  "country_id"),
    @("C8",  "customer_id"),
    @("C9",  "film_id"),
    @("C12", "inventory_id"),
    @("C13", "language_id"),
    @("C14", "payment_id"),
    @("C15", "rental_id"),
    @("C16", "staff_id"),
    @("C17", "store_id"),
    @("C10", "actor_id\film_id"),
    @("C11", "film_id\category_id")
)

foreach ($edit in $orderedEdits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

[void]$ws.Range("C11").Select()
